$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the content of C2 (previously contained "Some text")
$ws.Range("C2").ClearContents()

# Update the active selection to C2
$ws.Range("C2").Select()
